$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MMregeda_wrapper.c", "double d2;", "double d2=0;"),
    @("MMregeda_wrapper.c", "double c2;", "double c2=0;"),
    @("MMregeda_wrapper.c", "double B2;", "double B2=0;"),
    @("MMregeda_wrapper.c", "double A2;", "double A2=0;")
)

$startRow = 55
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

[void]$ws.Range("A55").Select()
